# Auto-generated edit script applying the Hades_Profits commit diff
# Updates computed market-price / profit figures across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1163.3334  # H18: 201 -> 1163.3334
$ws.Cells.Item(18, 9).Value = 1163.3334  # I18: 201 -> 1163.3334
$ws.Cells.Item(18, 11).Value = 1163.3334  # K18: 201 -> 1163.3334
$ws.Cells.Item(18, 13).Value = -879.3334  # M18: 83 -> -879.3334
$ws.Cells.Item(41, 8).Value = 342.33334  # H41: 334 -> 342.33334
$ws.Cells.Item(41, 9).Value = 470.2  # I41: 450.2 -> 470.2
$ws.Cells.Item(41, 11).Value = 470.2  # K41: 450.2 -> 470.2
$ws.Cells.Item(41, 13).Value = -30.19999999999999  # M41: -10.19999999999999 -> -30.19999999999999
$ws.Cells.Item(112, 8).Value = 40002092  # H112: 47621296 -> 40002092
$ws.Cells.Item(112, 10).Value = 2252.7273  # J112: 2468.889 -> 2252.7273
$ws.Cells.Item(112, 12).Value = 6758.1819  # L112: 7406.667 -> 6758.1819
$ws.Cells.Item(112, 14).Value = -8974.1819  # N112: -9622.667000000001 -> -8974.1819
$ws.Cells.Item(113, 8).Value = 3778.0454  # H113: 3663.28 -> 3778.0454
$ws.Cells.Item(113, 9).Value = 3702.6365  # I113: 3511.4614 -> 3702.6365
$ws.Cells.Item(113, 10).Value = 3853.4546  # J113: 3827.75 -> 3853.4546
$ws.Cells.Item(113, 11).Value = 3702.6365  # K113: 3511.4614 -> 3702.6365
$ws.Cells.Item(113, 12).Value = 3853.4546  # L113: 3827.75 -> 3853.4546
$ws.Cells.Item(113, 13).Value = -448.6365000000001  # M113: -257.4614000000001 -> -448.6365000000001
$ws.Cells.Item(113, 14).Value = -10361.4546  # N113: -10335.75 -> -10361.4546
$ws.Cells.Item(116, 8).Value = 1740.7333  # H116: 1934 -> 1740.7333
$ws.Cells.Item(116, 9).Value = 1593.2222  # I116: 1793 -> 1593.2222
$ws.Cells.Item(116, 10).Value = 1962  # J116: 2075 -> 1962
$ws.Cells.Item(116, 11).Value = 1593.2222  # K116: 1793 -> 1593.2222
$ws.Cells.Item(116, 12).Value = 1962  # L116: 2075 -> 1962
$ws.Cells.Item(116, 13).Value = 1848.7778  # M116: 1649 -> 1848.7778
$ws.Cells.Item(116, 14).Value = -8846  # N116: -8959 -> -8846
$ws.Cells.Item(132, 8).Value = 3009.8572  # H132: 2736.1729 -> 3009.8572
$ws.Cells.Item(132, 9).Value = 2822.6316  # I132: 2491.3635 -> 2822.6316
$ws.Cells.Item(132, 10).Value = 3830.7693  # J132: 3813.3333 -> 3830.7693
$ws.Cells.Item(132, 11).Value = 8467.8948  # K132: 7474.0905 -> 8467.8948
$ws.Cells.Item(132, 12).Value = 11492.3079  # L132: 11439.9999 -> 11492.3079
$ws.Cells.Item(132, 13).Value = -5937.8948  # M132: -4944.0905 -> -5937.8948
$ws.Cells.Item(132, 14).Value = -16552.3079  # N132: -16499.9999 -> -16552.3079
$ws.Cells.Item(135, 8).Value = 29754.139  # H135: 30666.885 -> 29754.139
$ws.Cells.Item(135, 9).Value = 36017.207  # I135: 37208.605 -> 36017.207
$ws.Cells.Item(135, 10).Value = 3807.1428  # J135: 4500 -> 3807.1428
$ws.Cells.Item(135, 11).Value = 324154.863  # K135: 334877.445 -> 324154.863
$ws.Cells.Item(135, 12).Value = 34264.2852  # L135: 40500 -> 34264.2852
$ws.Cells.Item(135, 13).Value = -321619.863  # M135: -332342.445 -> -321619.863
$ws.Cells.Item(135, 14).Value = -39334.2852  # N135: -45570 -> -39334.2852
$ws.Cells.Item(138, 8).Value = 2441070.2  # H138: 2085497.2 -> 2441070.2
$ws.Cells.Item(138, 9).Value = 1212.0233  # I138: 1295.641 -> 1212.0233
$ws.Cells.Item(138, 10).Value = 5131170.5  # J138: 3511529.8 -> 5131170.5
$ws.Cells.Item(138, 11).Value = 3636.0699  # K138: 3886.923 -> 3636.0699
$ws.Cells.Item(138, 12).Value = 15393511.5  # L138: 10534589.4 -> 15393511.5
$ws.Cells.Item(138, 13).Value = 1503.9301  # M138: 1253.077 -> 1503.9301
$ws.Cells.Item(138, 14).Value = -15403791.5  # N138: -10544869.4 -> -15403791.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 262.5  # H5: 323.16666 -> 262.5
$ws.Cells.Item(5, 9).Value = 187.14285  # I5: 229.8 -> 187.14285
$ws.Cells.Item(5, 11).Value = 187.14285  # K5: 229.8 -> 187.14285
$ws.Cells.Item(5, 13).Value = -75.14285000000001  # M5: -117.8 -> -75.14285000000001
$ws.Cells.Item(32, 8).Value = 4133.73  # H32: 3936.85 -> 4133.73
$ws.Cells.Item(32, 9).Value = 3953.5154  # I32: 3750.5464 -> 3953.5154
$ws.Cells.Item(32, 11).Value = 3953.5154  # K32: 3750.5464 -> 3953.5154
$ws.Cells.Item(32, 13).Value = -3666.5154  # M32: -3463.5464 -> -3666.5154
$ws.Cells.Item(61, 8).Value = 55668430  # H61: 45546892 -> 55668430
$ws.Cells.Item(61, 9).Value = 71501840  # I61: 50051480 -> 71501840
$ws.Cells.Item(61, 10).Value = 251498.5  # J61: 501000 -> 251498.5
$ws.Cells.Item(61, 11).Value = 71501840  # K61: 50051480 -> 71501840
$ws.Cells.Item(61, 12).Value = 251498.5  # L61: 501000 -> 251498.5
$ws.Cells.Item(61, 13).Value = -71501628  # M61: -50051268 -> -71501628
$ws.Cells.Item(61, 14).Value = -251922.5  # N61: -501424 -> -251922.5
$ws.Cells.Item(136, 8).Value = 55668430  # H136: 45546892 -> 55668430
$ws.Cells.Item(136, 9).Value = 71501840  # I136: 50051480 -> 71501840
$ws.Cells.Item(136, 10).Value = 251498.5  # J136: 501000 -> 251498.5
$ws.Cells.Item(136, 11).Value = 214505520  # K136: 150154440 -> 214505520
$ws.Cells.Item(136, 12).Value = 754495.5  # L136: 1503000 -> 754495.5
$ws.Cells.Item(136, 13).Value = -214502970  # M136: -150151890 -> -214502970
$ws.Cells.Item(136, 14).Value = -759595.5  # N136: -1508100 -> -759595.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 262.5  # H4: 323.16666 -> 262.5
$ws.Cells.Item(4, 9).Value = 187.14285  # I4: 229.8 -> 187.14285
$ws.Cells.Item(4, 11).Value = 187.14285  # K4: 229.8 -> 187.14285
$ws.Cells.Item(4, 13).Value = -72.14285000000001  # M4: -114.8 -> -72.14285000000001
$ws.Cells.Item(134, 8).Value = 25651916  # H134: 10422231 -> 25651916
$ws.Cells.Item(134, 9).Value = 13397.1  # I134: 6117.654 -> 13397.1
$ws.Cells.Item(134, 10).Value = 111113650  # J134: 55558724 -> 111113650
$ws.Cells.Item(134, 11).Value = 40191.3  # K134: 18352.962 -> 40191.3
$ws.Cells.Item(134, 12).Value = 333340950  # L134: 166676172 -> 333340950
$ws.Cells.Item(134, 13).Value = -37656.3  # M134: -15817.962 -> -37656.3
$ws.Cells.Item(134, 14).Value = -333346020  # N134: -166681242 -> -333346020

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 2083771.2  # H22: 2315234.8 -> 2083771.2
$ws.Cells.Item(22, 9).Value = 2604464  # I22: 2976444.5 -> 2604464
$ws.Cells.Item(22, 11).Value = 2604464  # K22: 2976444.5 -> 2604464
$ws.Cells.Item(22, 13).Value = -2604114  # M22: -2976094.5 -> -2604114
$ws.Cells.Item(31, 8).Value = 57629.824  # H31: 62227.945 -> 57629.824
$ws.Cells.Item(31, 9).Value = 69499.60000000001  # I31: 86644.586 -> 69499.60000000001
$ws.Cells.Item(31, 11).Value = 69499.60000000001  # K31: 86644.586 -> 69499.60000000001
$ws.Cells.Item(31, 13).Value = -69204.60000000001  # M31: -86349.586 -> -69204.60000000001
$ws.Cells.Item(34, 8).Value = 57629.824  # H34: 62227.945 -> 57629.824
$ws.Cells.Item(34, 9).Value = 69499.60000000001  # I34: 86644.586 -> 69499.60000000001
$ws.Cells.Item(34, 11).Value = 69499.60000000001  # K34: 86644.586 -> 69499.60000000001
$ws.Cells.Item(34, 13).Value = -69297.60000000001  # M34: -86442.586 -> -69297.60000000001
$ws.Cells.Item(107, 8).Value = 659.4138  # H107: 784.3 -> 659.4138
$ws.Cells.Item(107, 9).Value = 587.5  # I107: 702.75 -> 587.5
$ws.Cells.Item(107, 10).Value = 777.0909  # J107: 906.625 -> 777.0909
$ws.Cells.Item(107, 11).Value = 587.5  # K107: 702.75 -> 587.5
$ws.Cells.Item(107, 12).Value = 777.0909  # L107: 906.625 -> 777.0909
$ws.Cells.Item(107, 13).Value = 1332.5  # M107: 1217.25 -> 1332.5
$ws.Cells.Item(107, 14).Value = -4617.0909  # N107: -4746.625 -> -4617.0909
$ws.Cells.Item(122, 8).Value = 2550.6843  # H122: 2296.762 -> 2550.6843
$ws.Cells.Item(122, 9).Value = 2117  # I122: 1700.2667 -> 2117
$ws.Cells.Item(122, 10).Value = 2866.0908  # J122: 3788 -> 2866.0908
$ws.Cells.Item(122, 11).Value = 6351  # K122: 5100.800099999999 -> 6351
$ws.Cells.Item(122, 12).Value = 8598.2724  # L122: 11364 -> 8598.2724
$ws.Cells.Item(122, 13).Value = -3901  # M122: -2650.800099999999 -> -3901
$ws.Cells.Item(122, 14).Value = -13498.2724  # N122: -16264 -> -13498.2724
$ws.Cells.Item(124, 8).Value = 0  # H124: 35000 -> 0
$ws.Cells.Item(124, 10).Value = 0  # J124: 35000 -> 0
$ws.Cells.Item(124, 12).Value = 0  # L124: 35000 -> 0
$ws.Cells.Item(124, 14).ClearContents()  # N124: remove (was -39910)
$ws.Cells.Item(134, 8).Value = 34287.535  # H134: 101305 -> 34287.535
$ws.Cells.Item(134, 9).Value = 948.5925999999999  # I134: 1379.5 -> 948.5925999999999
$ws.Cells.Item(134, 10).Value = 334338  # J134: 501007 -> 334338
$ws.Cells.Item(134, 11).Value = 2845.7778  # K134: 4138.5 -> 2845.7778
$ws.Cells.Item(134, 12).Value = 1003014  # L134: 1503021 -> 1003014
$ws.Cells.Item(134, 13).Value = -310.7777999999998  # M134: -1603.5 -> -310.7777999999998
$ws.Cells.Item(134, 14).Value = -1008084  # N134: -1508091 -> -1008084

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 840.1429000000001  # H34: 863.1 -> 840.1429000000001
$ws.Cells.Item(34, 10).Value = 1502  # J34: 1209.3334 -> 1502
$ws.Cells.Item(34, 12).Value = 4506  # L34: 3628.0002 -> 4506
$ws.Cells.Item(34, 14).Value = -4674  # N34: -3796.0002 -> -4674
$ws.Cells.Item(39, 8).Value = 3548.5715  # H39: 4208 -> 3548.5715
$ws.Cells.Item(39, 10).Value = 4133.3335  # J39: 5250 -> 4133.3335
$ws.Cells.Item(39, 12).Value = 12400.0005  # L39: 15750 -> 12400.0005
$ws.Cells.Item(39, 14).Value = -12988.0005  # N39: -16338 -> -12988.0005
$ws.Cells.Item(112, 8).Value = 22225930  # H112: 30307758 -> 22225930
$ws.Cells.Item(112, 9).Value = 1500  # I112: 0 -> 1500
$ws.Cells.Item(112, 10).Value = 25645072  # J112: 30307758 -> 25645072
$ws.Cells.Item(112, 11).Value = 4500  # K112: 0 -> 4500
$ws.Cells.Item(112, 12).Value = 76935216  # L112: 90923274 -> 76935216
$ws.Cells.Item(112, 13).Value = -3392  # M112: None -> -3392
$ws.Cells.Item(112, 14).Value = -76937432  # N112: -90925490 -> -76937432
$ws.Cells.Item(131, 8).Value = 851.614  # H131: 854.322 -> 851.614
$ws.Cells.Item(131, 9).Value = 369.8889  # I131: 377.66666 -> 369.8889
$ws.Cells.Item(131, 10).Value = 941.9375  # J131: 940.12 -> 941.9375
$ws.Cells.Item(131, 11).Value = 1109.6667  # K131: 1132.99998 -> 1109.6667
$ws.Cells.Item(131, 12).Value = 2825.8125  # L131: 2820.36 -> 2825.8125
$ws.Cells.Item(131, 13).Value = 3930.3333  # M131: 3907.00002 -> 3930.3333
$ws.Cells.Item(131, 14).Value = -12905.8125  # N131: -12900.36 -> -12905.8125
$ws.Cells.Item(139, 8).Value = 4609.1  # H139: 4730.6895 -> 4609.1
$ws.Cells.Item(139, 10).Value = 7137.7334  # J139: 7570.2144 -> 7137.7334
$ws.Cells.Item(139, 12).Value = 21413.2002  # L139: 22710.6432 -> 21413.2002
$ws.Cells.Item(139, 14).Value = -31693.2002  # N139: -32990.6432 -> -31693.2002
$ws.Cells.Item(141, 8).Value = 15790  # H141: 14103.333 -> 15790
$ws.Cells.Item(141, 9).Value = 6046  # I141: 6832.5 -> 6046
$ws.Cells.Item(141, 10).Value = 19850  # J141: 19920 -> 19850
$ws.Cells.Item(141, 11).Value = 18138  # K141: 20497.5 -> 18138
$ws.Cells.Item(141, 12).Value = 59550  # L141: 59760 -> 59550
$ws.Cells.Item(141, 13).Value = -12958  # M141: -15317.5 -> -12958
$ws.Cells.Item(141, 14).Value = -69910  # N141: -70120 -> -69910

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 1815.8572  # H126: 1781.8889 -> 1815.8572
$ws.Cells.Item(126, 9).Value = 1742.2  # I126: 1646.1428 -> 1742.2
$ws.Cells.Item(126, 10).Value = 2000  # J126: 2257 -> 2000
$ws.Cells.Item(126, 11).Value = 5226.6  # K126: 4938.428400000001 -> 5226.6
$ws.Cells.Item(126, 12).Value = 6000  # L126: 6771 -> 6000
$ws.Cells.Item(126, 13).Value = -2756.6  # M126: -2468.428400000001 -> -2756.6
$ws.Cells.Item(126, 14).Value = -10940  # N126: -11711 -> -10940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 765.8125  # H22: 756.2 -> 765.8125
$ws.Cells.Item(22, 9).Value = 744.1177  # I22: 730.55 -> 744.1177
$ws.Cells.Item(22, 11).Value = 744.1177  # K22: 730.55 -> 744.1177
$ws.Cells.Item(22, 13).Value = -449.1177  # M22: -435.55 -> -449.1177
$ws.Cells.Item(27, 8).Value = 765.8125  # H27: 756.2 -> 765.8125
$ws.Cells.Item(27, 9).Value = 744.1177  # I27: 730.55 -> 744.1177
$ws.Cells.Item(27, 11).Value = 744.1177  # K27: 730.55 -> 744.1177
$ws.Cells.Item(27, 13).Value = -637.1177  # M27: -623.55 -> -637.1177
$ws.Cells.Item(55, 8).Value = 208.82353  # H55: 240.6875 -> 208.82353
$ws.Cells.Item(55, 9).Value = 229.08333  # I55: 245.75 -> 229.08333
$ws.Cells.Item(55, 10).Value = 160.2  # J55: 225.5 -> 160.2
$ws.Cells.Item(55, 11).Value = 229.08333  # K55: 245.75 -> 229.08333
$ws.Cells.Item(55, 12).Value = 160.2  # L55: 225.5 -> 160.2
$ws.Cells.Item(55, 13).Value = -56.08332999999999  # M55: -72.75 -> -56.08332999999999
$ws.Cells.Item(55, 14).Value = -506.2  # N55: -571.5 -> -506.2
$ws.Cells.Item(58, 8).Value = 4500  # H58: 0 -> 4500
$ws.Cells.Item(58, 9).Value = 4500  # I58: 0 -> 4500
$ws.Cells.Item(58, 11).Value = 4500  # K58: 0 -> 4500
$ws.Cells.Item(58, 13).Value = -4240  # M58: None -> -4240
$ws.Cells.Item(61, 8).Value = 2469.5715  # H61: 2758.6 -> 2469.5715
$ws.Cells.Item(61, 9).Value = 2463.7334  # I61: 2819.4 -> 2463.7334
$ws.Cells.Item(61, 10).Value = 2484.1667  # J61: 2637 -> 2484.1667
$ws.Cells.Item(61, 11).Value = 2463.7334  # K61: 2819.4 -> 2463.7334
$ws.Cells.Item(61, 12).Value = 2484.1667  # L61: 2637 -> 2484.1667
$ws.Cells.Item(61, 13).Value = -2261.7334  # M61: -2617.4 -> -2261.7334
$ws.Cells.Item(61, 14).Value = -2888.1667  # N61: -3041 -> -2888.1667
$ws.Cells.Item(82, 8).Value = 1887.4  # H82: 2107.9167 -> 1887.4
$ws.Cells.Item(82, 9).Value = 640.5  # I82: 663.3333 -> 640.5
$ws.Cells.Item(82, 10).Value = 2718.6667  # J82: 2589.4443 -> 2718.6667
$ws.Cells.Item(82, 11).Value = 640.5  # K82: 663.3333 -> 640.5
$ws.Cells.Item(82, 12).Value = 2718.6667  # L82: 2589.4443 -> 2718.6667
$ws.Cells.Item(82, 13).Value = -279.5  # M82: -302.3333 -> -279.5
$ws.Cells.Item(82, 14).Value = -3440.6667  # N82: -3311.4443 -> -3440.6667
$ws.Cells.Item(85, 8).Value = 1887.4  # H85: 2107.9167 -> 1887.4
$ws.Cells.Item(85, 9).Value = 640.5  # I85: 663.3333 -> 640.5
$ws.Cells.Item(85, 10).Value = 2718.6667  # J85: 2589.4443 -> 2718.6667
$ws.Cells.Item(85, 11).Value = 640.5  # K85: 663.3333 -> 640.5
$ws.Cells.Item(85, 12).Value = 2718.6667  # L85: 2589.4443 -> 2718.6667
$ws.Cells.Item(85, 13).Value = 607.5  # M85: 584.6667 -> 607.5
$ws.Cells.Item(85, 14).Value = -5214.6667  # N85: -5085.4443 -> -5214.6667
$ws.Cells.Item(93, 8).Value = 1016.2647  # H93: 973.9167 -> 1016.2647
$ws.Cells.Item(93, 9).Value = 972.3333  # I93: 922.7931 -> 972.3333
$ws.Cells.Item(93, 11).Value = 972.3333  # K93: 922.7931 -> 972.3333
$ws.Cells.Item(93, 13).Value = 275.6667  # M93: 325.2069 -> 275.6667
$ws.Cells.Item(100, 8).Value = 1550.409  # H100: 1595.5238 -> 1550.409
$ws.Cells.Item(100, 9).Value = 1100.8182  # I100: 1150.6 -> 1100.8182
$ws.Cells.Item(100, 11).Value = 1100.8182  # K100: 1150.6 -> 1100.8182
$ws.Cells.Item(100, 13).Value = -559.8181999999999  # M100: -609.5999999999999 -> -559.8181999999999
$ws.Cells.Item(113, 8).Value = 2469.5715  # H113: 2758.6 -> 2469.5715
$ws.Cells.Item(113, 9).Value = 2463.7334  # I113: 2819.4 -> 2463.7334
$ws.Cells.Item(113, 10).Value = 2484.1667  # J113: 2637 -> 2484.1667
$ws.Cells.Item(113, 11).Value = 2463.7334  # K113: 2819.4 -> 2463.7334
$ws.Cells.Item(113, 12).Value = 2484.1667  # L113: 2637 -> 2484.1667
$ws.Cells.Item(113, 13).Value = -293.7334000000001  # M113: -649.4000000000001 -> -293.7334000000001
$ws.Cells.Item(113, 14).Value = -6824.1667  # N113: -6977 -> -6824.1667
$ws.Cells.Item(132, 8).Value = 22413.29  # H132: 27583.12 -> 22413.29
$ws.Cells.Item(132, 9).Value = 2532.875  # I132: 3175.0833 -> 2532.875
$ws.Cells.Item(132, 10).Value = 54221.95  # J132: 60127.168 -> 54221.95
$ws.Cells.Item(132, 11).Value = 7598.625  # K132: 9525.249899999999 -> 7598.625
$ws.Cells.Item(132, 12).Value = 162665.85  # L132: 180381.504 -> 162665.85
$ws.Cells.Item(132, 13).Value = -5068.625  # M132: -6995.249899999999 -> -5068.625
$ws.Cells.Item(132, 14).Value = -167725.85  # N132: -185441.504 -> -167725.85

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 40666.668  # H46: 0 -> 40666.668
$ws.Cells.Item(46, 10).Value = 40666.668  # J46: 0 -> 40666.668
$ws.Cells.Item(46, 12).Value = 40666.668  # L46: 0 -> 40666.668
$ws.Cells.Item(46, 14).Value = -41128.668  # N46: None -> -41128.668
$ws.Cells.Item(122, 8).Value = 1861.3125  # H122: 1631.7222 -> 1861.3125
$ws.Cells.Item(122, 9).Value = 1361.5  # I122: 1207.4546 -> 1361.5
$ws.Cells.Item(122, 10).Value = 2361.125  # J122: 2298.4285 -> 2361.125
$ws.Cells.Item(122, 11).Value = 4084.5  # K122: 3622.3638 -> 4084.5
$ws.Cells.Item(122, 12).Value = 7083.375  # L122: 6895.2855 -> 7083.375
$ws.Cells.Item(122, 13).Value = -1634.5  # M122: -1172.3638 -> -1634.5
$ws.Cells.Item(122, 14).Value = -11983.375  # N122: -11795.2855 -> -11983.375
$ws.Cells.Item(132, 8).Value = 29048.57  # H132: 34229.23 -> 29048.57
$ws.Cells.Item(132, 9).Value = 27205.13  # I132: 31248.031 -> 27205.13
$ws.Cells.Item(132, 10).Value = 31108.883  # J132: 37742.785 -> 31108.883
$ws.Cells.Item(132, 11).Value = 81615.39  # K132: 93744.09299999999 -> 81615.39
$ws.Cells.Item(132, 12).Value = 93326.649  # L132: 113228.355 -> 93326.649
$ws.Cells.Item(132, 13).Value = -79085.39  # M132: -91214.09299999999 -> -79085.39
$ws.Cells.Item(132, 14).Value = -98386.649  # N132: -118288.355 -> -98386.649
$ws.Cells.Item(134, 8).Value = 40666.668  # H134: 0 -> 40666.668
$ws.Cells.Item(134, 10).Value = 40666.668  # J134: 0 -> 40666.668
$ws.Cells.Item(134, 12).Value = 122000.004  # L134: 0 -> 122000.004
$ws.Cells.Item(134, 14).Value = -127070.004  # N134: None -> -127070.004
$ws.Cells.Item(136, 8).Value = 31271.324  # H136: 37210.754 -> 31271.324
$ws.Cells.Item(136, 9).Value = 20256.46  # I136: 24388.441 -> 20256.46
$ws.Cells.Item(136, 10).Value = 67069.625  # J136: 76593.57000000001 -> 67069.625
$ws.Cells.Item(136, 11).Value = 60769.38  # K136: 73165.323 -> 60769.38
$ws.Cells.Item(136, 12).Value = 201208.875  # L136: 229780.71 -> 201208.875
$ws.Cells.Item(136, 13).Value = -58219.38  # M136: -70615.323 -> -58219.38
$ws.Cells.Item(136, 14).Value = -206308.875  # N136: -234880.71 -> -206308.875
